# Atualizei dados bibi e add
# Update row 9 (Ano 2025) of the faturamento_anual sheet with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3406867.06
$ws.Range("C9").Value = 532306.99
$ws.Range("D9").Value = 3939174.05
$ws.Range("E9").Value = 13.51316248643545
$ws.Range("F9").Value = 86.48683751356455
$ws.Range("G9").Value = -48.5551836215154
$ws.Range("H9").Value = -38.47669000301178
$ws.Range("I9").Value = 34115
$ws.Range("J9").Value = 1445
$ws.Range("K9").Value = 35560
$ws.Range("L9").Value = 24537
$ws.Range("M9").Value = 160.5401658719485
$ws.Range("N9").Value = 9.604041079437131
